$d = $word.ActiveDocument

$p = $d.Paragraphs.Item(30)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="HTMLconformatoprevio"/><w:shd w:val="clear" w:color="auto" w:fill="1E1F22"/><w:ind w:left="720"/><w:rPr><w:color w:val="A9B7C6"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="E8BF6A"/></w:rPr><w:lastRenderedPageBreak/><w:t>&lt;uses-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="E8BF6A"/></w:rPr><w:t>permission</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="E8BF6A"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="9876AA"/></w:rPr><w:t>android</w:t></w:r><w:r><w:rPr><w:color w:val="BABABA"/></w:rPr><w:t>:name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="6A8759"/></w:rPr><w:t>="</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="6A8759"/></w:rPr><w:t>android.permission.Internet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="6A8759"/></w:rPr><w:t>"</w:t></w:r><w:r><w:rPr><w:color w:val="E8BF6A"/></w:rPr><w:t>&gt;&lt;/uses-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="E8BF6A"/></w:rPr><w:t>permission</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="E8BF6A"/></w:rPr><w:t>&gt;</w:t></w:r></w:p>')

$p = $d.Paragraphs.Item(29)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Se habilita los permisos de Internet en el archivo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Manifest</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p = $d.Paragraphs.Item(27)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Se crea el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fragment</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> denominado: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>registrarProductoFragment</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y luego crear todos los campos necesarios</w:t></w:r></w:p>')

$p = $d.Paragraphs.Item(23)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="HTMLconformatoprevio"/><w:shd w:val="clear" w:color="auto" w:fill="1E1F22"/><w:ind w:left="360"/><w:rPr><w:color w:val="A9B7C6"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="9876AA"/></w:rPr><w:t>android</w:t></w:r><w:r><w:rPr><w:color w:val="BABABA"/></w:rPr><w:t>:usesCleartextTraffic</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="6A8759"/></w:rPr><w:t>="true"</w:t></w:r></w:p>')

$p = $d.Paragraphs.Item(22)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">La librería </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Volley</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> por defecto solo permite hacer peticiones en HTTPS, se recomienda tener un servidor con este protocolo, en caso que no sea posible, se puede hacer el siguiente ajuste en el archivo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>manifest</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, para que acepte este tipo de petición.</w:t></w:r></w:p>')

$p = $d.Paragraphs.Item(20)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:shd w:val="clear" w:color="auto" w:fill="1E1F22"/><w:tabs><w:tab w:val="left" w:pos="916"/><w:tab w:val="left" w:pos="1832"/><w:tab w:val="left" w:pos="2748"/><w:tab w:val="left" w:pos="3664"/><w:tab w:val="left" w:pos="4580"/><w:tab w:val="left" w:pos="5496"/><w:tab w:val="left" w:pos="6412"/><w:tab w:val="left" w:pos="7328"/><w:tab w:val="left" w:pos="8244"/><w:tab w:val="left" w:pos="9160"/><w:tab w:val="left" w:pos="10076"/><w:tab w:val="left" w:pos="10992"/><w:tab w:val="left" w:pos="11908"/><w:tab w:val="left" w:pos="12824"/><w:tab w:val="left" w:pos="13740"/><w:tab w:val="left" w:pos="14656"/></w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="A9B7C6"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CO"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="808080"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CO"/><w14:ligatures w14:val="none"/></w:rPr><w:t>//se implementa la siguiente librería</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="808080"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CO"/><w14:ligatures w14:val="none"/></w:rPr><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:i/><w:iCs/><w:color w:val="FFC66D"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CO"/><w14:ligatures w14:val="none"/></w:rPr><w:t>implementation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="A9B7C6"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CO"/><w14:ligatures w14:val="none"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="6A8759"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CO"/><w14:ligatures w14:val="none"/></w:rPr><w:t>"com.android.volley:volley:1.2.1"</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="A9B7C6"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CO"/><w14:ligatures w14:val="none"/></w:rPr><w:t>)</w:t></w:r></w:p>')

$p = $d.Paragraphs.Item(19)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">En el archivo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>build.gradle</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">   (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Module:app</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>) la siguiente librería</w:t></w:r></w:p>')

$p = $d.Paragraphs.Item(18)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Crear proyecto denominado: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>CRUDMovil</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p = $d.Paragraphs.Item(16)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>JDK</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p = $d.Paragraphs.Item(13)
$p.Range.InsertXML('<w:p><w:pPr><w:ind w:left="2832"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>function</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>eliminarProducto</w:t></w:r><w:r><w:t>&amp;id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=valor1</w:t></w:r></w:p>')

$p = $d.Paragraphs.Item(10)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>consultar</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Producto</w:t></w:r><w:r><w:t xml:space="preserve"> por ID</w:t></w:r><w:r><w:t xml:space="preserve"> con los siguientes parámetros:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="2880"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>function</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>consultarProductoPorCodigo</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="2880"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>guardar</w:t></w:r><w:r><w:t xml:space="preserve"> p</w:t></w:r><w:r><w:t>roducto</w:t></w:r><w:r><w:t xml:space="preserve"> con los siguientes parámetros: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>function</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>guardarProducto&amp;id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=</w:t></w:r><w:r><w:t>valor1</w:t></w:r><w:r><w:t>&amp;nombre=</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>valor2</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>&amp;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>descripcion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>valor3</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>&amp;precio=</w:t></w:r><w:r><w:t>valor4</w:t></w:r><w:r><w:t>&amp;cantidad=</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>valor5</w:t></w:r><w:r><w:t>&amp;imagen=</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>valor6</w:t></w:r></w:p>')

$p = $d.Paragraphs.Item(8)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="2880"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>function</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>consultarListaProductos</w:t></w:r><w:r><w:t>&amp;filtro</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=valor</w:t></w:r></w:p>')

$p = $d.Paragraphs.Item(7)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="2880"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>function</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>consultarListaProductos</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p = $d.Paragraphs.Item(4)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:ind w:left="2124"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>controller</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>productosController</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p = $d.Paragraphs.Item(2)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Se crea un back-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>end</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">  con la siguiente información</w:t></w:r></w:p>')
